# Updates "想去人数" (want-to-go count) figures that were refreshed by the
# site's data-generation job (gh-pages output regenerated at commit 456a3b4).
# Same event rows appear both on the "展览" sheet and on the aggregated
# "全部类型" sheet, so both need to be kept in sync.

$wb = $excel.ActiveWorkbook

# -- Sheet "展览" --------------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 7655
$ws1.Range("F3").Value  = 7655
$ws1.Range("F6").Value  = 38
$ws1.Range("F9").Value  = 6586
$ws1.Range("F10").Value = 3356
$ws1.Range("F20").Value = 22
$ws1.Range("F21").Value = 310
$ws1.Range("F23").Value = 3817
$ws1.Range("F28").Value = 1464
$ws1.Range("F34").Value = 43
$ws1.Range("F42").Value = 1406

# -- Sheet "全部类型" ------------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value  = 7655
$ws4.Range("F6").Value  = 7655
$ws4.Range("F11").Value = 6586
$ws4.Range("F12").Value = 3356
$ws4.Range("F21").Value = 310
$ws4.Range("F24").Value = 3817
$ws4.Range("F32").Value = 1464
$ws4.Range("F38").Value = 43
$ws4.Range("F46").Value = 1406
